# Generate Report for Handoff
#
# The localization-status report gains a freshly generated "Latest Handoff
# Datetime" for the c7ddd372-29fd-45ba-9639-119b8bad870d row (row 7) on both
# the "zh-cn" and "de-de" per-language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-08 14:30:50"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-08 14:30:55"
